$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 255.77777
$ws.Range("I33").Value = 118.2
$ws.Range("J33").Value = 427.75
$ws.Range("K33").Value = 118.2
$ws.Range("L33").Value = 427.75
$ws.Range("M33").Value = 110.8
$ws.Range("N33").Value = -885.75

$ws.Range("H62").Value = 348939.56
$ws.Range("I62").Value = 479548.47
$ws.Range("J62").Value = 6091.125
$ws.Range("K62").Value = 479548.47
$ws.Range("L62").Value = 6091.125
$ws.Range("M62").Value = -478924.47
$ws.Range("N62").Value = -7339.125

$ws.Range("H65").Value = 348939.56
$ws.Range("I65").Value = 479548.47
$ws.Range("J65").Value = 6091.125
$ws.Range("K65").Value = 2397742.35
$ws.Range("L65").Value = 30455.625
$ws.Range("M65").Value = -2394622.35
$ws.Range("N65").Value = -36695.625

$ws.Range("H106").Value = 1679
$ws.Range("I106").Value = 1679
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1679
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1048
$ws.Range("N106").ClearContents()

$ws.Range("H112").Value = 1706
$ws.Range("I112").Value = 778
$ws.Range("J112").Value = 2015.3334
$ws.Range("K112").Value = 2334
$ws.Range("L112").Value = 6046.0002
$ws.Range("M112").Value = -1226
$ws.Range("N112").Value = -8262.0002

$ws.Range("H125").Value = 967.4375
$ws.Range("I125").Value = 961.44446
$ws.Range("K125").Value = 8653.00014
$ws.Range("M125").Value = -6193.00014

$ws.Range("H138").Value = 3139.75
$ws.Range("I138").Value = 1060.7368
$ws.Range("J138").Value = 4413.984
$ws.Range("K138").Value = 3182.2104
$ws.Range("L138").Value = 13241.952
$ws.Range("M138").Value = 1957.7896
$ws.Range("N138").Value = -23521.952

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13324.13
$ws.Range("I32").Value = 7563.5806
$ws.Range("J32").Value = 89857.14
$ws.Range("K32").Value = 7563.5806
$ws.Range("L32").Value = 89857.14
$ws.Range("M32").Value = -7276.5806
$ws.Range("N32").Value = -90431.14

$ws.Range("H80").Value = 26262.5
$ws.Range("J80").Value = 26262.5
$ws.Range("L80").Value = 26262.5
$ws.Range("N80").Value = -28258.5

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H83").Value = 26262.5
$ws.Range("J83").Value = 26262.5
$ws.Range("L83").Value = 78787.5
$ws.Range("N83").Value = -88771.5

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H132").Value = 1647.1111
$ws.Range("I132").Value = 1189.1034
$ws.Range("J132").Value = 3544.5715
$ws.Range("K132").Value = 3567.3102
$ws.Range("L132").Value = 10633.7145
$ws.Range("M132").Value = -1037.3102
$ws.Range("N132").Value = -15693.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 756.7646999999999
$ws.Range("I94").Value = 795.21277
$ws.Range("J94").Value = 305
$ws.Range("K94").Value = 795.21277
$ws.Range("L94").Value = 305
$ws.Range("M94").Value = -344.21277
$ws.Range("N94").Value = -1207

$ws.Range("H105").Value = 9153.333000000001
$ws.Range("I105").Value = 9153.333000000001
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 9153.333000000001
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -7406.333000000001
$ws.Range("N105").ClearContents()

$ws.Range("H134").Value = 1846.0465
$ws.Range("I134").Value = 1600
$ws.Range("K134").Value = 4800
$ws.Range("M134").Value = -2265

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1388.8276
$ws.Range("I132").Value = 1134
$ws.Range("J132").Value = 2189.7144
$ws.Range("K132").Value = 3402
$ws.Range("L132").Value = 6569.1432
$ws.Range("M132").Value = -872
$ws.Range("N132").Value = -11629.1432

$ws.Range("H134").Value = 23810808
$ws.Range("I134").Value = 1348.25
$ws.Range("J134").Value = 500000000
$ws.Range("K134").Value = 4044.75
$ws.Range("L134").Value = 1500000000
$ws.Range("M134").Value = -1509.75
$ws.Range("N134").Value = -1500005070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3676.9
$ws.Range("I3").Value = 844.8333
$ws.Range("K3").Value = 2534.4999
$ws.Range("M3").Value = -2422.4999

$ws.Range("H46").Value = 9333.333000000001
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 9333.333000000001
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 27999.999
$ws.Range("N46").Value = -28181.999
$ws.Range("M46").ClearContents()

$ws.Range("H86").Value = 122
$ws.Range("I86").Value = 122
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 366
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 820
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 122
$ws.Range("I89").Value = 122
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1098
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 4830
$ws.Range("N89").ClearContents()

$ws.Range("H113").Value = 563.46
$ws.Range("I113").Value = 794
$ws.Range("J113").Value = 498.43588
$ws.Range("K113").Value = 2382
$ws.Range("L113").Value = 1495.30764
$ws.Range("M113").Value = -212
$ws.Range("N113").Value = -5835.30764

$ws.Range("H136").Value = 74082.78999999999
$ws.Range("I136").Value = 112462.11
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 337386.33
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -332286.33
$ws.Range("N136").Value = -25200

$ws.Range("H137").Value = 4523.857
$ws.Range("I137").Value = 3051.3333
$ws.Range("J137").Value = 8205.166999999999
$ws.Range("K137").Value = 9153.999899999999
$ws.Range("L137").Value = 24615.501
$ws.Range("M137").Value = -4053.999899999999
$ws.Range("N137").Value = -34815.501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2067.2727
$ws.Range("I97").Value = 2090.7693
$ws.Range("J97").Value = 2033.3334
$ws.Range("K97").Value = 2090.7693
$ws.Range("L97").Value = 2033.3334
$ws.Range("M97").Value = -1594.7693
$ws.Range("N97").Value = -3025.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1370.7059
$ws.Range("I100").Value = 1240.1538
$ws.Range("J100").Value = 1795
$ws.Range("K100").Value = 2480.3076
$ws.Range("L100").Value = 3590
$ws.Range("M100").Value = -1939.3076
$ws.Range("N100").Value = -4672

$ws.Range("H122").Value = 41667988
$ws.Range("I122").Value = 47620370
$ws.Range("J122").Value = 1333.3334
$ws.Range("K122").Value = 142861110
$ws.Range("L122").Value = 4000.0002
$ws.Range("M122").Value = -142858660
$ws.Range("N122").Value = -8900.0002
